# Update the "District" column (G) on the sheet so that district names use
# the current official names, and drop a few stray empty "Address" (F)
# cells that shouldn't have held a value in the first place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-55 (data rows) all get their District (column G) value normalized
# to "Belagavi (Belgaum)", except:
#   - row 24, whose G cell held a stray duplicate of the person's name and
#     should read the (unrelated) official district name "Vijayapura (Bijapur)"
#   - row 49, whose G cell is left untouched
for ($r = 3; $r -le 55; $r++) {
    if ($r -eq 49) { continue }
    if ($r -eq 24) {
        $ws.Cells.Item($r, 7).Value = "Vijayapura (Bijapur)"
    } else {
        $ws.Cells.Item($r, 7).Value = "Belagavi (Belgaum)"
    }
}

# These four rows had a spurious empty inline-string cell in column F
# (Address) that carries no data - remove it entirely.
$ws.Range("F24").ClearContents()
$ws.Range("F25").ClearContents()
$ws.Range("F34").ClearContents()
$ws.Range("F52").ClearContents()
